# Updates the cryptos price/volume table (Sheet1) to the latest scrape.
# Source: "Updated cryptos list on Mon Mar  4 17:57:01 UTC 2024 with GitHub Actions"
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '67.146.00'
$ws.Range("E2").Value = '  +7.25%  '

$ws.Range("D3").Value = '3.554.14'
$ws.Range("E3").Value = '  +2.61%  '

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.999'
$ws.Range("E4").Value = '  -0.18%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '415.35'
$ws.Range("E5").Value = '  +0.17%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '128.35'
$ws.Range("E6").Value = '  -1.57%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.644'
$ws.Range("E7").Value = '  +3.34%  '

$ws.Range("D8").Value = '3.540.95'
$ws.Range("E8").Value = '  +2.40%  '

$ws.Range("E9").Value = '  -0.02%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.772'
$ws.Range("E10").Value = '  +5.85%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.176'
$ws.Range("E11").Value = '  +25.67%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.0000324'
$ws.Range("E12").Value = '  +47.51%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '42.21'
$ws.Range("E13").Value = '  -1.28%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '9.83'
$ws.Range("E14").Value = '  +1.64%  '

$ws.Range("D15").Value = '4.133.55'
$ws.Range("E15").Value = '  +2.87%  '

$ws.Range("E16").Value = '  -0.21%  '

$ws.Range("B17").Value = 'WrappedEther'
$ws.Range("C17").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D17").Value = '3.642.75'
$ws.Range("E17").Value = '  +4.85%  '

$ws.Range("B18").Value = 'Chainlink'
$ws.Range("C18").Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '20.03'
$ws.Range("E18").Value = '  -2.26%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '1.11'
$ws.Range("E19").Value = '  +2.58%  '

$ws.Range("D20").Value = '66.890.77'
$ws.Range("E20").Value = '  +6.83%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '12.23'
$ws.Range("E21").Value = '  -4.36%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '458.66'
$ws.Range("E22").Value = '  -2.70%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '89.25'
$ws.Range("E23").Value = '  -1.63%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '3.13'
$ws.Range("E24").Value = '  -4.30%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '12.86'
$ws.Range("E25").Value = '  -4.12%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '3.43'
$ws.Range("E26").Value = '  +3.60%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '9.80'
$ws.Range("E27").Value = '  -7.43%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '34.77'
$ws.Range("E28").Value = '  +4.22%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '4.86'
$ws.Range("E29").Value = '  +1.22%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '2.75'
$ws.Range("E30").Value = '  +4.04%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '12.20'
$ws.Range("E31").Value = '  +1.37%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.115'
$ws.Range("E32").Value = '  +2.22%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '7.20'
$ws.Range("E33").Value = '  -5.35%  '

$ws.Range("E34").Value = '  -6.31%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.997'
$ws.Range("E35").Value = '  -0.22%  '

$ws.Range("E36").Value = '  -6.02%  '

$ws.Range("E37").Value = '  -4.01%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.0488'
$ws.Range("E38").Value = '  -0.46%  '

$ws.Range("D39").Value = '0.0₃0747'
$ws.Range("E39").Value = '  +30.78%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.146'
$ws.Range("E40").Value = '  +8.78%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '1.00'
$ws.Range("E41").Value = '  +0.10%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '147.23'
$ws.Range("E42").Value = '  +1.57%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '2.94'
$ws.Range("E43").Value = '  -3.78%  '

$ws.Range("E44").Value = '  +0.26%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '4.28'
$ws.Range("E45").Value = '  -2.05%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '3.18'
$ws.Range("E46").Value = '  -5.22%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.302'
$ws.Range("E47").Value = '  -6.08%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '1.94'
$ws.Range("E48").Value = '  -6.39%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '119.41'
$ws.Range("E49").Value = '  +8.89%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '2.25'
$ws.Range("E50").Value = '  -5.75%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '2.57'
$ws.Range("E51").Value = '  +9.55%  '
